$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Adm/Calcr -> sCs)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.650723999999999
$ws.Range("H2").Value = 25.952172
$ws.Range("I2").Value = 0.2392268437287548
$ws.Range("J2").Value = 0.2392268437287548
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.194209333333333
$ws.Range("N2").Value = 15.582628
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 44.93367134089065
$ws.Range("R2").Value = 404.403042068016
$ws.Range("S2").Value = 0.2392268437287548
$ws.Range("T2").Value = 0.2392268437287548

# Row 3 (FAPs -> Adm/Calcr -> sCs)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Adm"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 26.269711
$ws.Range("H3").Value = 78.809133
$ws.Range("I3").Value = 0.7264617444963627
$ws.Range("J3").Value = 0.7264617444963627
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.194209333333333
$ws.Range("N3").Value = 15.582628
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 136.4503780601693
$ws.Range("R3").Value = 1228.053402541524
$ws.Range("S3").Value = 0.7264617444963627
$ws.Range("T3").Value = 0.7264617444963627

# Row 4 (M1 -> Adm/Calcr -> sCs) - NEW ROW
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Adm"
$ws.Range("C4").Value = "Calcr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.021265
$ws.Range("H4").Value = 0.063795
$ws.Range("I4").Value = 0.0005880616271992926
$ws.Range("J4").Value = 0.0005880616271992925
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.194209333333333
$ws.Range("N4").Value = 15.582628
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.1104548614733333
$ws.Range("R4").Value = 0.99409375326
$ws.Range("S4").Value = 0.0005880616271992926
$ws.Range("T4").Value = 0.0005880616271992925

# Row 5 (Neutro -> Adm/Calcr -> sCs) - NEW ROW
$ws.Range("A5").Value = "Neutro"
$ws.Range("B5").Value = "Adm"
$ws.Range("C5").Value = "Calcr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3188806666666666
$ws.Range("H5").Value = 0.956642
$ws.Range("I5").Value = 0.008818315717018348
$ws.Range("J5").Value = 0.008818315717018348
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.194209333333333
$ws.Range("N5").Value = 15.582628
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1.656332935019555
$ws.Range("R5").Value = 14.906996415176
$ws.Range("S5").Value = 0.008818315717018348
$ws.Range("T5").Value = 0.008818315717018348

# Row 6 (sCs -> Adm/Calcr -> sCs) - NEW ROW
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Adm"
$ws.Range("C6").Value = "Calcr"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9005953333333334
$ws.Range("H6").Value = 2.701786
$ws.Range("I6").Value = 0.02490503443066491
$ws.Range("J6").Value = 0.02490503443066491
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.194209333333333
$ws.Range("N6").Value = 15.582628
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 4.677880685956445
$ws.Range("R6").Value = 42.100926173608
$ws.Range("S6").Value = 0.02490503443066491
$ws.Range("T6").Value = 0.02490503443066491
